$wb = $excel.ActiveWorkbook

# Rename sheets
$wsInclude = $wb.Worksheets.Item("Include from PayorToPayorExch")
$wsInclude.Name = "Include #0"

$wsValueSets = $wb.Worksheets.Item("Include ValueSets")
$wsValueSets.Name = "Include ValueSet #1"

# Update the Metadata sheet
$ws = $wb.Worksheets.Item("Metadata")

# Version: 0.2.1 -> 0.2.2
$ws.Range("B3").Value = "0.2.2"

# Date: updated
$ws.Range("B8").Value = "2024-09-11T16:17:59-05:00"

# Contact: updated
$ws.Range("B10").Value = "MITRE (https://www.mitre.org)"

# Insert a new "Jurisdiction" row right after the Contact row (row 10),
# pushing Description/Purpose/Copyright/Immutable down by one.
$ws.Rows.Item(11).Insert()

$ws.Range("A11").Value = "Jurisdiction"
# Leading apostrophe forces a literal (empty) text value rather than
# clearing the cell, so B11 ends up as an explicit empty string.
$ws.Range("B11").Value = "'"

# Copy formatting from the row now below (originally row 11, the
# "Description" row) so the new row matches the sheet's normal cell style.
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
